$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# The sheet had two duplicate "Contact" / "No display for ContactDetail" rows
# (rows 10 and 11). Remove the second one so everything below shifts up by one
# row, then retarget the remaining row to the new Jurisdiction property.
$meta.Rows.Item(11).Delete()

$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Sheet "Elements" (sheet2) ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element; Short/Definition now describe the
# specific last-claim-indicator extension instead of the generic text.
$elements.Range("K2").Value = "Last Claim Indicator"
$elements.Range("L2").Value = "Indicates whether this claim record is the last or most recent claim"
